$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.799.84'
$ws.Range('D3').Value = '1.896.89'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''0.7649'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').Value = '''240.14'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.3034'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').Value = '''25.33'
$ws.Range('E9').Value = '  -3.83%  '
$ws.Range('D10').Value = '''0.06819'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = '''0.07971'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '''0.7335'
$ws.Range('E12').Value = '  -4.64%  '
$ws.Range('D13').Value = '1.868.90'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').Value = '''5.153'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '29.798.33'
$ws.Range('D17').Value = '''13.77'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '''5.887'
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').Value = '''241.54'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = '''0.000007688'
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('D21').Value = '''0.9998'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '''1.001'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '''6.878'
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '''166.33'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '''9.202'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''18.60'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '''0.1284'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''2.021'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''1.400'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.512'
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''4.249'
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''4.057'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.05198'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''1.244'
$ws.Range('E34').Value = '  -2.23%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.7229'
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.714'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.01912'
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.774'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''6.142'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.4382'
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''71.70'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.001'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '''1.879'
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '''0.8272'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''7.588'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '''99.51'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''9.690'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.027.55'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''36.04'
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05921'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').Value = '''1.469'
$ws.Range('E51').Value = '  +1.69%  '
